$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.908.90"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.815.01"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.49"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4655"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07347"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8688"
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.818.89"
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.384"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07082"
$ws.Range("E14").Value = "  +0.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.511"
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.64"
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008707"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.66"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.944.43"
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.074.49"
$ws.Range("E24").Value = "  +0.86%  "
$ws.Range("E25").Value = "  -0.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.97"
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.38"
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.142"
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.255"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.56"
$ws.Range("E30").Value = "  -0.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08920"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7562"
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.485"
$ws.Range("E34").Value = "  +0.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.912"
$ws.Range("E35").Value = "  -0.98%  "
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("E37").Value = "  -1.39%  "
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.978"
$ws.Range("E40").Value = "  +1.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.248"
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5303"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.286"
$ws.Range("E43").Value = "  -3.50%  "
$ws.Range("E44").Value = "  -0.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.432"
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4875"
$ws.Range("E46").Value = "  -2.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.38"
$ws.Range("E47").Value = "  +0.93%  "
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("E49").Value = "  -0.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.660"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06292"
$ws.Range("E51").Value = "  +0.10%  "
